$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 543.3570999999999
$ws.Range("I28").Value = 444.45456
$ws.Range("J28").Value = 906
$ws.Range("K28").Value = 444.45456
$ws.Range("L28").Value = 906
$ws.Range("M28").Value = 40.54543999999999
$ws.Range("N28").Value = -1876
$ws.Range("H41").Value = 651.1579
$ws.Range("J41").Value = 915.2
$ws.Range("L41").Value = 915.2
$ws.Range("N41").Value = -1795.2
$ws.Range("H62").Value = 200004770
$ws.Range("J62").Value = 7457.5
$ws.Range("L62").Value = 7457.5
$ws.Range("N62").Value = -8705.5
$ws.Range("H65").Value = 200004770
$ws.Range("J65").Value = 7457.5
$ws.Range("L65").Value = 37287.5
$ws.Range("N65").Value = -43527.5
$ws.Range("H96").Value = 801.5454999999999
$ws.Range("I96").Value = 461.8
$ws.Range("J96").Value = 1084.6666
$ws.Range("K96").Value = 1385.4
$ws.Range("L96").Value = 3253.9998
$ws.Range("M96").Value = -12.40000000000009
$ws.Range("N96").Value = -5999.9998
$ws.Range("H138").Value = 144204.66
$ws.Range("J138").Value = 4994.338
$ws.Range("L138").Value = 14983.014
$ws.Range("N138").Value = -25263.014

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 698890.2
$ws.Range("I32").Value = 698890.2
$ws.Range("K32").Value = 698890.2
$ws.Range("M32").Value = -698603.2
$ws.Range("H45").Value = 64899
$ws.Range("I45").Value = 83949.92
$ws.Range("J45").Value = 5364.875
$ws.Range("K45").Value = 83949.92
$ws.Range("L45").Value = 5364.875
$ws.Range("M45").Value = -83572.92
$ws.Range("N45").Value = -6118.875
$ws.Range("H46").Value = 10721.8
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 10902.25
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 10902.25
$ws.Range("M46").Value = -9681
$ws.Range("N46").Value = -11540.25
$ws.Range("H61").Value = 5039.7095
$ws.Range("I61").Value = 4742.579
$ws.Range("K61").Value = 4742.579
$ws.Range("M61").Value = -4530.579
$ws.Range("H63").Value = 1090.5
$ws.Range("I63").Value = 1088
$ws.Range("J63").Value = 1093
$ws.Range("K63").Value = 1088
$ws.Range("L63").Value = 1093
$ws.Range("M63").Value = -402
$ws.Range("N63").Value = -2465
$ws.Range("H66").Value = 1090.5
$ws.Range("I66").Value = 1088
$ws.Range("J66").Value = 1093
$ws.Range("K66").Value = 5440
$ws.Range("L66").Value = 5465
$ws.Range("M66").Value = -2008
$ws.Range("N66").Value = -12329
$ws.Range("H102").Value = 11176.934
$ws.Range("I102").Value = 13016.211
$ws.Range("J102").Value = 8000
$ws.Range("K102").Value = 13016.211
$ws.Range("L102").Value = 8000
$ws.Range("M102").Value = -11394.211
$ws.Range("N102").Value = -11244
$ws.Range("H132").Value = 2850.8
$ws.Range("I132").Value = 1949
$ws.Range("J132").Value = 5813.857
$ws.Range("K132").Value = 5847
$ws.Range("L132").Value = 17441.571
$ws.Range("M132").Value = -3317
$ws.Range("N132").Value = -22501.571
$ws.Range("H136").Value = 5039.7095
$ws.Range("I136").Value = 4742.579
$ws.Range("K136").Value = 14227.737
$ws.Range("M136").Value = -11677.737

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 827.2
$ws.Range("I22").Value = 382.33334
$ws.Range("K22").Value = 382.33334
$ws.Range("M22").Value = -209.33334
$ws.Range("H64").Value = 5908.2666
$ws.Range("I64").Value = 15629.8
$ws.Range("J64").Value = 1047.5
$ws.Range("K64").Value = 15629.8
$ws.Range("L64").Value = 1047.5
$ws.Range("M64").Value = -15404.8
$ws.Range("N64").Value = -1497.5
$ws.Range("H67").Value = 5908.2666
$ws.Range("I67").Value = 15629.8
$ws.Range("J67").Value = 1047.5
$ws.Range("K67").Value = 15629.8
$ws.Range("L67").Value = 1047.5
$ws.Range("M67").Value = -14849.8
$ws.Range("N67").Value = -2607.5
$ws.Range("H82").Value = 23916
$ws.Range("I82").Value = 8101
$ws.Range("J82").Value = 59499.75
$ws.Range("K82").Value = 8101
$ws.Range("L82").Value = 59499.75
$ws.Range("M82").Value = -7718
$ws.Range("N82").Value = -60265.75
$ws.Range("H85").Value = 23916
$ws.Range("I85").Value = 8101
$ws.Range("J85").Value = 59499.75
$ws.Range("K85").Value = 8101
$ws.Range("L85").Value = 59499.75
$ws.Range("M85").Value = -6775
$ws.Range("N85").Value = -62151.75
$ws.Range("H99").Value = 39941.9
$ws.Range("I99").Value = 70486
$ws.Range("K99").Value = 70486
$ws.Range("M99").Value = -68988
$ws.Range("H132").Value = 77199
$ws.Range("J132").Value = 77199
$ws.Range("L132").Value = 77199
$ws.Range("N132").Value = -87319
$ws.Range("H134").Value = 1926.6904
$ws.Range("I134").Value = 1555.8422
$ws.Range("K134").Value = 4667.5266
$ws.Range("M134").Value = -2132.5266

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3208.1943
$ws.Range("I31").Value = 2287.0435
$ws.Range("J31").Value = 4837.923
$ws.Range("K31").Value = 2287.0435
$ws.Range("L31").Value = 4837.923
$ws.Range("M31").Value = -1992.0435
$ws.Range("N31").Value = -5427.923
$ws.Range("H34").Value = 3208.1943
$ws.Range("I34").Value = 2287.0435
$ws.Range("J34").Value = 4837.923
$ws.Range("K34").Value = 2287.0435
$ws.Range("L34").Value = 4837.923
$ws.Range("M34").Value = -2085.0435
$ws.Range("N34").Value = -5241.923
$ws.Range("H58").Value = 2317.575
$ws.Range("J58").Value = 2447.4666
$ws.Range("L58").Value = 2447.4666
$ws.Range("N58").Value = -2853.4666
$ws.Range("H99").Value = 3750200.2
$ws.Range("I99").Value = 6454400.5
$ws.Range("J99").Value = 5922.769
$ws.Range("K99").Value = 6454400.5
$ws.Range("L99").Value = 5922.769
$ws.Range("M99").Value = -6452902.5
$ws.Range("N99").Value = -8918.769
$ws.Range("H122").Value = 6954.0386
$ws.Range("I122").Value = 8772.588
$ws.Range("K122").Value = 26317.764
$ws.Range("M122").Value = -23867.764
$ws.Range("H125").Value = 75000
$ws.Range("J125").Value = 75000
$ws.Range("L125").Value = 75000
$ws.Range("N125").Value = -79920
$ws.Range("H126").Value = 3750200.2
$ws.Range("I126").Value = 6454400.5
$ws.Range("J126").Value = 5922.769
$ws.Range("K126").Value = 19363201.5
$ws.Range("L126").Value = 17768.307
$ws.Range("M126").Value = -19360731.5
$ws.Range("N126").Value = -22708.307
$ws.Range("H132").Value = 6333.1816
$ws.Range("I132").Value = 6746.5
$ws.Range("K132").Value = 20239.5
$ws.Range("M132").Value = -17709.5
$ws.Range("H134").Value = 3626.7144
$ws.Range("I134").Value = 3606.1667
$ws.Range("K134").Value = 10818.5001
$ws.Range("M134").Value = -8283.500100000001
$ws.Range("H136").Value = 2317.575
$ws.Range("J136").Value = 2447.4666
$ws.Range("L136").Value = 7342.399800000001
$ws.Range("N136").Value = -12442.3998
$ws.Range("H141").Value = 397378.84
$ws.Range("J141").Value = 446962.1
$ws.Range("L141").Value = 446962.1
$ws.Range("N141").Value = -457322.1

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 41671360
$ws.Range("I131").Value = 76930056
$ws.Range("J131").Value = 1989.7273
$ws.Range("K131").Value = 230790168
$ws.Range("L131").Value = 5969.1819
$ws.Range("M131").Value = -230785128
$ws.Range("N131").Value = -16049.1819

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 105.92308
$ws.Range("I2").Value = 111.7
$ws.Range("J2").Value = 86.666664
$ws.Range("K2").Value = 111.7
$ws.Range("L2").Value = 86.666664
$ws.Range("M2").Value = 1.299999999999997
$ws.Range("N2").Value = -312.666664
$ws.Range("H4").Value = 3666
$ws.Range("I4").Value = 4998
$ws.Range("K4").Value = 4998
$ws.Range("M4").Value = -4886
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 2812.58
$ws.Range("I132").Value = 2509.8293
$ws.Range("K132").Value = 7529.4879
$ws.Range("M132").Value = -4999.4879

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H122").Value = 9582.916999999999
$ws.Range("I122").Value = 13748.75
$ws.Range("K122").Value = 41246.25
$ws.Range("M122").Value = -38796.25
$ws.Range("H132").Value = 406761.2
$ws.Range("I132").Value = 1151364.4
$ws.Range("K132").Value = 3454093.2
$ws.Range("M132").Value = -3451563.2

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1273.5294
$ws.Range("I4").Value = 116.666664
$ws.Range("J4").Value = 1521.4286
$ws.Range("K4").Value = 116.666664
$ws.Range("L4").Value = 1521.4286
$ws.Range("M4").Value = -3.666663999999997
$ws.Range("N4").Value = -1747.4286
$ws.Range("H44").Value = 24064
$ws.Range("J44").Value = 24064
$ws.Range("L44").Value = 24064
$ws.Range("N44").Value = -25172
$ws.Range("H107").Value = 15633.762
$ws.Range("I107").Value = 1480.4667
$ws.Range("J107").Value = 51017
$ws.Range("K107").Value = 4441.4001
$ws.Range("L107").Value = 153051
$ws.Range("M107").Value = -2521.4001
$ws.Range("N107").Value = -156891
$ws.Range("H132").Value = 3876.9055
$ws.Range("J132").Value = 3676.1
$ws.Range("L132").Value = 11028.3
$ws.Range("N132").Value = -16088.3
$ws.Range("H136").Value = 331699.3
$ws.Range("I136").Value = 362082.97
$ws.Range("J136").Value = 5075
$ws.Range("K136").Value = 1086248.91
$ws.Range("L136").Value = 15225
$ws.Range("M136").Value = -1083698.91
$ws.Range("N136").Value = -20325
